$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'68.465.14"
$ws.Range("E2").Value = '  -4.21%  '

$ws.Range("D3").Value = "'3.684.61"
$ws.Range("E3").Value = '  -5.19%  '

$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = '  +0.03%  '

$ws.Range("D5").Value = "'593.89"
$ws.Range("E5").Value = '  -2.17%  '

$ws.Range("D6").Value = "'181.68"
$ws.Range("E6").Value = '  +5.15%  '

$ws.Range("D7").Value = "'3.681.09"
$ws.Range("E7").Value = '  -5.12%  '

$ws.Range("D8").Value = "'0.631"
$ws.Range("E8").Value = '  -6.13%  '

$ws.Range("D9").Value = "'0.999"
$ws.Range("E9").Value = '  -0.04%  '

$ws.Range("D10").Value = "'0.715"
$ws.Range("E10").Value = '  -5.15%  '

$ws.Range("E11").Value = '  -10.12%  '

$ws.Range("D12").Value = "'55.91"
$ws.Range("E12").Value = '  +3.28%  '

$ws.Range("D13").Value = "'0.0000293"
$ws.Range("E13").Value = '  -9.73%  '

$ws.Range("D14").Value = "'10.68"
$ws.Range("E14").Value = '  -7.67%  '

$ws.Range("D15").Value = "'4.262.72"

$ws.Range("D16").Value = "'3.678.00"
$ws.Range("E16").Value = '  -5.18%  '

$ws.Range("D17").Value = "'19.38"
$ws.Range("E17").Value = '  -8.88%  '

$ws.Range("E18").Value = '  -2.26%  '

$ws.Range("B19").Value = 'Uniswap'
$ws.Range("C19").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D19").Value = "'12.84"
$ws.Range("E19").Value = '  -8.32%  '

$ws.Range("B20").Value = 'Polygon'
$ws.Range("C20").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D20").Value = "'1.12"
$ws.Range("E20").Value = '  -7.74%  '

$ws.Range("D21").Value = "'68.248.76"
$ws.Range("E21").Value = '  -4.21%  '

$ws.Range("D22").Value = "'410.96"
$ws.Range("E22").Value = '  -7.10%  '

$ws.Range("D23").Value = "'4.59"
$ws.Range("E23").Value = '  -4.78%  '

$ws.Range("D24").Value = "'88.73"
$ws.Range("E24").Value = '  -6.68%  '

$ws.Range("D25").Value = "'3.03"
$ws.Range("E25").Value = '  -8.75%  '

$ws.Range("D26").Value = "'12.68"
$ws.Range("E26").Value = '  -9.25%  '

$ws.Range("D27").Value = "'10.77"
$ws.Range("E27").Value = '  -9.40%  '

$ws.Range("D28").Value = "'3.90"
$ws.Range("E28").Value = '  -3.75%  '

$ws.Range("D29").Value = "'6.08"
$ws.Range("E29").Value = '  +1.86%  '

$ws.Range("D30").Value = "'9.60"
$ws.Range("E30").Value = '  -9.03%  '

$ws.Range("D31").Value = "'32.85"
$ws.Range("E31").Value = '  -7.11%  '

$ws.Range("D32").Value = "'7.39"
$ws.Range("E32").Value = '  -16.08%  '

$ws.Range("D33").Value = "'12.44"
$ws.Range("E33").Value = '  -8.65%  '

$ws.Range("D34").Value = "'0.120"
$ws.Range("E34").Value = '  -5.82%  '

$ws.Range("D35").Value = "'64.69"
$ws.Range("E35").Value = '  -6.86%  '

$ws.Range("D36").Value = "'43.40"
$ws.Range("E36").Value = '  -9.82%  '

$ws.Range("D37").Value = "'604.52"
$ws.Range("E37").Value = '  -5.22%  '

$ws.Range("D38").Value = "'0.0₃0884"
$ws.Range("E38").Value = '  -11.66%  '

$ws.Range("E39").Value = '  +0.05%  '

$ws.Range("D40").Value = "'0.402"
$ws.Range("E40").Value = '  -9.52%  '

$ws.Range("E41").Value = '  -0.06%  '

$ws.Range("D42").Value = "'0.138"
$ws.Range("E42").Value = '  -6.77%  '

$ws.Range("D43").Value = "'3.03"
$ws.Range("E43").Value = '  -7.76%  '

$ws.Range("B44").Value = 'VeChain'
$ws.Range("C44").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D44").Value = "'0.0442"
$ws.Range("E44").Value = '  -6.89%  '

$ws.Range("B45").Value = 'Fetch.AI'
$ws.Range("C45").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D45").Value = "'2.65"
$ws.Range("E45").Value = '  -8.84%  '

$ws.Range("D46").Value = "'2.75"
$ws.Range("E46").Value = '  -14.03%  '

$ws.Range("B47").Value = 'WEMIXToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D47").Value = "'2.73"
$ws.Range("E47").Value = '  -6.30%  '

$ws.Range("B48").Value = 'Stellar'
$ws.Range("C48").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D48").Value = "'0.135"
$ws.Range("E48").Value = '  -6.52%  '

$ws.Range("D49").Value = "'9.00"
$ws.Range("E49").Value = '  -12.01%  '

$ws.Range("D50").Value = "'2.720.64"
$ws.Range("E50").Value = '  -7.19%  '

$ws.Range("E51").Value = '  -5.32%  '
